$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row translation (Spanish -> English snake_case) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case "de/del/la/las/el/los/y" connectors in state/municipality names (and 1 typo fix) ---
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San José De Gracia"
$ws.Range("B15").Value = "Playas De Rosarito"
$ws.Range("B36").Value = "Amatenango De La Frontera"
$ws.Range("B37").Value = "Amatenango Del Valle"
$ws.Range("B40").Value = "Bejucal De Ocampo"
$ws.Range("B42").Value = "Benemérito De Las Américas"
$ws.Range("B51").Value = "Chiapa De Corzo"
$ws.Range("B57").Value = "Comitán De Domínguez"
$ws.Range("B85").Value = "Marqués De Comillas"
$ws.Range("B86").Value = "Mazapa De Madero"
$ws.Range("B90").Value = "Montecristo De Guerrero"
$ws.Range("B94").Value = "Ocozocoautla De Espinosa"
$ws.Range("B106").Value = "Salto De Agua"
$ws.Range("B108").Value = "San Cristóbal De Las Casas"
$ws.Range("B112").Value = "Santiago El Pinar"
$ws.Range("B149").Value = "Coyame Del Sotol"
$ws.Range("B154").Value = "Guadalupe Y Calvo"
$ws.Range("B156").Value = "Hidalgo Del Parral"
$ws.Range("B189").Value = "San Juan De Sabinas"
$ws.Range("B199").Value = "Villa De Álvarez"
$ws.Range("A201").Value = "Ciudad De México"
$ws.Range("B205").Value = "Cuajimalpa De Morelos"
$ws.Range("B230").Value = "Nombre De Dios"
$ws.Range("B234").Value = "Pánuco De Coronado"
$ws.Range("B241").Value = "San Juan De Guadalupe"
$ws.Range("B242").Value = "San Juan Del Río"
$ws.Range("A251").Value = "Estado De México"
$ws.Range("B251").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B254").Value = "Almoloya De Alquisiras"
$ws.Range("B255").Value = "Almoloya De Juárez"
$ws.Range("B260").Value = "Atizapán De Zaragoza"
$ws.Range("B265").Value = "Chapa De Mota"
$ws.Range("B269").Value = "Coacalco De Berriozábal"
$ws.Range("B275").Value = "Ecatepec De Morelos"
$ws.Range("B281").Value = "Ixtapan De La Sal"
$ws.Range("B295").Value = "Naucalpan De Juárez"
$ws.Range("B307").Value = "San Felipe Del Progreso"
$ws.Range("B309").Value = "San Simón De Guerrero"
$ws.Range("B319").Value = "Tenango Del Valle"
$ws.Range("B329").Value = "Tlalnepantla De Baz"
$ws.Range("B335").Value = "Valle De Bravo"
$ws.Range("B336").Value = "Valle De Chalco Solidaridad"
$ws.Range("B337").Value = "Villa De Allende"
$ws.Range("B338").Value = "Villa Del Carbón"
$ws.Range("B352").Value = "Apaseo El Alto"
$ws.Range("B353").Value = "Apaseo El Grande"
$ws.Range("B361").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B365").Value = "Jaral Del Progreso"
$ws.Range("B373").Value = "Purísima Del Rincón"
$ws.Range("B377").Value = "San Diego De La Unión"
$ws.Range("B379").Value = "San Francisco Del Rincón"
$ws.Range("B381").Value = "San Luis De La Paz"
$ws.Range("B383").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B385").Value = "Silao De La Victoria"
$ws.Range("B390").Value = "Valle De Santiago"
$ws.Range("B396").Value = "Acapulco De Juárez"
$ws.Range("B399").Value = "Ajuchitlán Del Progreso"
$ws.Range("B400").Value = "Alcozauca De Guerrero"
$ws.Range("B404").Value = "Atenango Del Río"
$ws.Range("B406").Value = "Atoyac De Álvarez"
$ws.Range("B407").Value = "Ayutla De Los Libres"
$ws.Range("B409").Value = "Buenavista De Cuéllar"
$ws.Range("B410").Value = "Chilapa De Álvarez"
$ws.Range("B411").Value = "Chilpancingo De Los Bravo"
$ws.Range("B412").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B417").Value = "Coyuca De Benítez"
$ws.Range("B418").Value = "Coyuca De Catalán"
$ws.Range("B422").Value = "Cuetzala Del Progreso"
$ws.Range("B423").Value = "Cutzamala De Pinzón"
$ws.Range("B429").Value = "Huitzuco De Los Figueroa"
$ws.Range("B430").Value = "Iguala De La Independencia"
$ws.Range("B432").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B433").Value = "Zihuatanejo De Azueta"
$ws.Range("B435").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B438").Value = "Mártir De Cuilapan"
$ws.Range("B451").Value = "Taxco De Alarcón"
$ws.Range("B453").Value = "Técpan De Galeana"
$ws.Range("B455").Value = "Tepecoacuilco De Trujano"
$ws.Range("B457").Value = "Tixtla De Guerrero"
$ws.Range("B461").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B462").Value = "Tlapa De Comonfort"
$ws.Range("B474").Value = "Agua Blanca De Iturbide"
$ws.Range("B480").Value = "Atotonilco El Grande"
$ws.Range("B486").Value = "Cuautepec De Hinojosa"
$ws.Range("B491").Value = "Huasca De Ocampo"
$ws.Range("B495").Value = "Huejutla De Reyes"
$ws.Range("B498").Value = "Jacala De Ledezma"
$ws.Range("B505").Value = "Mineral De La Reforma"
$ws.Range("B506").Value = "Mineral Del Chico"
$ws.Range("B507").Value = "Mineral Del Monte"
$ws.Range("B508").Value = "Mixquiahuala De Juárez"
$ws.Range("B509").Value = "Molango De Escamilla"
$ws.Range("B511").Value = "Nopala De Villagrán"
$ws.Range("B512").Value = "Omitlán De Juárez"
$ws.Range("B513").Value = "Pachuca De Soto"
$ws.Range("B516").Value = "Progreso De Obregón"
$ws.Range("B522").Value = "Santiago De Anaya"
$ws.Range("B523").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B527").Value = "Tenango De Doria"
$ws.Range("B529").Value = "Tepehuacán De Guerrero"
$ws.Range("B530").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B532").Value = "Tezontepec De Aldama"
$ws.Range("B540").Value = "Tula De Allende"
$ws.Range("B541").Value = "Tulancingo De Bravo"
$ws.Range("B542").Value = "Villa De Tezontepec"
$ws.Range("B545").Value = "Zacualtipán De Ángeles"
$ws.Range("B554").Value = "Atotonilco El Alto"
$ws.Range("B555").Value = "Autlán De Navarro"
$ws.Range("B561").Value = "Cañadas De Obregón"
$ws.Range("B566").Value = "Cuautitlán De García Barragán"
$ws.Range("B573").Value = "Encarnación De Díaz"
$ws.Range("B579").Value = "Huejuquilla El Alto"
$ws.Range("B580").Value = "Ixtlahuacán Del Río"
$ws.Range("B584").Value = "Jilotlán De Los Dolores"
$ws.Range("B589").Value = "Lagos De Moreno"
$ws.Range("B595").Value = "Ojuelos De Jalisco"
$ws.Range("B600").Value = "San Cristóbal De La Barranca"
$ws.Range("B601").Value = "San Diego De Alejandría"
$ws.Range("B603").Value = "San Martín De Bolaños"
$ws.Range("B605").Value = "San Miguel El Alto"
$ws.Range("B606").Value = "San Sebastián Del Oeste"
$ws.Range("B607").Value = "Santa María De Los Ángeles"
$ws.Range("B608").Value = "Santa María Del Oro"
$ws.Range("B611").Value = "Talpa De Allende"
$ws.Range("B612").Value = "Tamazula De Gordiano"
$ws.Range("B616").Value = "Teocuitatlán De Corona"
$ws.Range("B617").Value = "Tepatitlán De Morelos"
$ws.Range("B619").Value = "Tizapán El Alto"
$ws.Range("B620").Value = "Tlajomulco De Zúñiga"
$ws.Range("B628").Value = "Unión De San Antonio"
$ws.Range("B629").Value = "Valle De Juárez"
$ws.Range("B633").Value = "Yahualica De González Gallo"
$ws.Range("B634").Value = "Zacoalco De Torres"
$ws.Range("B637").Value = "Zapotlán Del Rey"
$ws.Range("B638").Value = "Zapotlán El Grande"
$ws.Range("B662").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B664").Value = "Cojumatlán De Régules"
$ws.Range("B728").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B752").Value = "Coatlán Del Río"
$ws.Range("B760").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B764").Value = "Puente De Ixtla"
$ws.Range("B770").Value = "Tetela Del Volcán"
$ws.Range("B772").Value = "Tlaltizapán De Zapata"
$ws.Range("B778").Value = "Zacualpan De Amilpas"
$ws.Range("B782").Value = "Bahía De Banderas"
$ws.Range("B785").Value = "Ixtlán Del Río"
$ws.Range("B791").Value = "Santa María Del Oro"
$ws.Range("B815").Value = "Mier Y Noriega"
$ws.Range("B816").Value = "Montemorelos"
$ws.Range("B819").Value = "San Nicolás De Los Garza"
$ws.Range("B825").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B833").Value = "Chalcatongo De Hidalgo"
$ws.Range("B834").Value = "Ciénega De Zimatlán"
$ws.Range("B837").Value = "Coicoyán De Las Flores"
$ws.Range("B840").Value = "Cuilápam De Guerrero"
$ws.Range("B841").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B842").Value = "Guadalupe De Ramírez"
$ws.Range("B843").Value = "Guevea De Humboldt"
$ws.Range("B844").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B845").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B846").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B847").Value = "Huautla De Jiménez"
$ws.Range("B849").Value = "Ixtlán De Juárez"
$ws.Range("B850").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B859").Value = "Mariscala De Juárez"
$ws.Range("B861").Value = "Mazatlán Villa De Flores"
$ws.Range("B863").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B864").Value = "Mixistlán De La Reforma"
$ws.Range("B868").Value = "Nejapa De Madero"
$ws.Range("B870").Value = "Oaxaca De Juárez"
$ws.Range("B871").Value = "Ocotlán De Morelos"
$ws.Range("B872").Value = "Pinotepa De Don Luis"
$ws.Range("B874").Value = "Putla Villa De Guerrero"
$ws.Range("B875").Value = "Reforma De Pineda"
$ws.Range("B890").Value = "San Antonino El Alto"
$ws.Range("B892").Value = "San Antonio De La Cal"
$ws.Range("B908").Value = "San Dionisio Del Mar"
$ws.Range("B911").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B916").Value = "San Francisco Del Mar"
$ws.Range("B935").Value = "San José Del Progreso"
$ws.Range("B942").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B950").Value = "San Juan Del Estado"
$ws.Range("B983").Value = "San Mateo Del Mar"
$ws.Range("B993").Value = "San Miguel Del Puerto"
$ws.Range("B994").Value = "San Miguel Del Río"
$ws.Range("B995").Value = "San Miguel El Grande"
$ws.Range("B1007").Value = "San Pablo Villa De Mitla"
$ws.Range("B1011").Value = "San Pedro El Alto"
$ws.Range("B1026").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B1027").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B1037").Value = "Santa Ana Del Valle"
$ws.Range("B1044").Value = "Santa Cruz De Bravo"
$ws.Range("B1052").Value = "Santa Lucía Del Camino"
$ws.Range("B1068").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B1093").Value = "Santiago Del Río"
$ws.Range("B1125").Value = "Santo Domingo De Morelos"
$ws.Range("B1145").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1146").Value = "Tanetze De Zaragoza"
$ws.Range("B1147").Value = "Tataltepec De Valdés"
$ws.Range("B1148").Value = "Teotitlán De Flores Magón"
$ws.Range("B1149").Value = "Teotitlán Del Valle"
$ws.Range("B1151").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B1152").Value = "Tlacolula De Matamoros"
$ws.Range("B1153").Value = "Tlalixtac De Cabrera"
$ws.Range("B1157").Value = "Villa De Chilapa De Díaz"
$ws.Range("B1158").Value = "Villa De Etla"
$ws.Range("B1159").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B1160").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B1161").Value = "Villa De Zaachila"
$ws.Range("B1162").Value = "Villa Sola De Vega"
$ws.Range("B1163").Value = "Zapotitlán Del Río"
$ws.Range("B1165").Value = "Zimatlán De Álvarez"
$ws.Range("B1183").Value = "Ayotoxco De Guerrero"
$ws.Range("B1186").Value = "Chalchicomula De Sesma"
$ws.Range("B1195").Value = "Chila De La Sal"
$ws.Range("B1203").Value = "Cuapiaxtla De Madero"
$ws.Range("B1207").Value = "Cuayuca De Andrade"
$ws.Range("B1208").Value = "Cuetzalan Del Progreso"
$ws.Range("B1222").Value = "Huehuetlán El Chico"
$ws.Range("B1223").Value = "Huehuetlán El Grande"
$ws.Range("B1228").Value = "Ixcamilpa De Guerrero"
$ws.Range("B1231").Value = "Izúcar De Matamoros"
$ws.Range("B1240").Value = "Los Reyes De Juárez"
$ws.Range("B1250").Value = "Palmar De Bravo"
$ws.Range("B1266").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1270").Value = "San Salvador El Seco"
$ws.Range("B1271").Value = "San Salvador El Verde"
$ws.Range("B1282").Value = "Tepanco De López"
$ws.Range("B1287").Value = "Tepexi De Rodríguez"
$ws.Range("B1289").Value = "Tetela De Ocampo"
$ws.Range("B1290").Value = "Teteles De Avila Castillo"
$ws.Range("B1295").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1306").Value = "Tuzamapan De Galeana"
$ws.Range("B1309").Value = "Xayacatlán De Bravo"
$ws.Range("B1326").Value = "Amealco De Bonfil"
$ws.Range("B1328").Value = "Cadereyta De Montes"
$ws.Range("B1334").Value = "Jalpan De Serra"
$ws.Range("B1335").Value = "Landa De Matamoros"
$ws.Range("B1338").Value = "Pinal De Amoles"
$ws.Range("B1341").Value = "San Juan Del Río"
$ws.Range("B1352").Value = "Armadillo De Los Infante"
$ws.Range("B1353").Value = "Axtla De Terrazas"
$ws.Range("B1358").Value = "Ciudad Del Maíz"
$ws.Range("B1368").Value = "Mexquitic De Carmona"
$ws.Range("B1373").Value = "San Ciro De Acosta"
$ws.Range("B1379").Value = "Santa María Del Río"
$ws.Range("B1381").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1387").Value = "Tanquián De Escobedo"
$ws.Range("B1391").Value = "Villa De Arista"
$ws.Range("B1392").Value = "Villa De Arriaga"
$ws.Range("B1393").Value = "Villa De Guadalupe"
$ws.Range("B1394").Value = "Villa De La Paz"
$ws.Range("B1395").Value = "Villa De Ramos"
$ws.Range("B1396").Value = "Villa De Reyes"
$ws.Range("B1446").Value = "Jalpa De Méndez"
$ws.Range("B1485").Value = "Soto La Marina"
$ws.Range("B1497").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1502").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1504").Value = "Muñoz De Domingo Arenas"
$ws.Range("B1506").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1509").Value = "San Pablo Del Monte"
$ws.Range("B1528").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1532").Value = "Amatlán De Los Reyes"
$ws.Range("B1544").Value = "Boca Del Río"
$ws.Range("B1548").Value = "Castillo De Teayo"
$ws.Range("B1550").Value = "Cazones De Herrera"
$ws.Range("B1557").Value = "Chinampa De Gorostiza"
$ws.Range("B1567").Value = "Cosamaloapan De Carpio"
$ws.Range("B1568").Value = "Cosautlán De Carvajal"
$ws.Range("B1584").Value = "Hueyapan De Ocampo"
$ws.Range("B1585").Value = "Huiloapan De Cuauhtémoc"
$ws.Range("B1586").Value = "Ignacio De La Llave"
$ws.Range("B1589").Value = "Ixhuatlán De Madero"
$ws.Range("B1590").Value = "Ixhuatlán Del Café"
$ws.Range("B1591").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1601").Value = "Juchique De Ferrer"
$ws.Range("B1605").Value = "Lerdo De Tejada"
$ws.Range("B1609").Value = "Martínez De La Torre"
$ws.Range("B1611").Value = "Medellín De Bravo"
$ws.Range("B1615").Value = "Mixtla De Altamirano"
$ws.Range("B1617").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1625").Value = "Ozuluama De Mascareñas"
$ws.Range("B1629").Value = "Paso Del Macho"
$ws.Range("B1632").Value = "Poza Rica De Hidalgo"
$ws.Range("B1641").Value = "Sayula De Alemán"
$ws.Range("B1644").Value = "Soledad De Doblado"
$ws.Range("B1651").Value = "Tatahuicapan De Juárez"
$ws.Range("B1685").Value = "Vega De Alatorre"
$ws.Range("B1694").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1695").Value = "Zozocolco De Hidalgo"
$ws.Range("B1712").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1714").Value = "Concepción Del Oro"
$ws.Range("B1716").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1725").Value = "Jiménez Del Teul"
$ws.Range("B1732").Value = "Nochistlán De Mejía"
$ws.Range("B1741").Value = "Teúl De González Ortega"
$ws.Range("B1742").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1746").Value = "Villa De Cos"

# --- Floating point precision touch-ups (27/27187 ratio cells) ---
$ws.Range("D300").Value = 0.0009931217125832201
$ws.Range("D380").Value = 0.0009931217125832201
$ws.Range("D430").Value = 0.0009931217125832201
$ws.Range("D524").Value = 0.0009931217125832201
$ws.Range("D687").Value = 0.0009931217125832201
$ws.Range("D763").Value = 0.0009931217125832201
$ws.Range("D1562").Value = 0.0009931217125832201
$ws.Range("D1624").Value = 0.0009931217125832201

# --- Remove trailing footnote rows (1755-1759) ---
$ws.Range("A1755:D1759").EntireRow.Delete()

Write-Host "Edit complete"